$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("F16").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F33").Value = -10
$ws.Range("F40").Value = -2
$ws.Range("F46").Value = 2
$ws.Range("F48").Value = 3
$ws.Range("F51").Value = -1
$ws.Range("F53").Value = 1
$ws.Range("F54").Value = 0
$ws.Range("F60").Value = 0
$ws.Range("F65").Value = 0
